$d = $word.ActiveDocument

$replacements = @(
    @{old="53×80=4240"; new="94×79=7426"},
    @{old="78×90=7020"; new="57×23=1311"},
    @{old="87×88=7656"; new="88×62=5456"},
    @{old="51×93=4743"; new="97×54=5238"},
    @{old="16×78=1248"; new="51×84=4284"},
    @{old="53×53=2809"; new="67×81=5427"},
    @{old="82×21=1722"; new="57×48=2736"},
    @{old="98×48=4704"; new="73×76=5548"},
    @{old="57×82=4674"; new="95×86=8170"},
    @{old="36×26=936"; new="28×40=1120"},
    @{old="67×37=2479"; new="27×59=1593"},
    @{old="52×34=1768"; new="19×58=1102"},
    @{old="53×35=1855"; new="93×31=2883"},
    @{old="50×61=3050"; new="30×18=540"},
    @{old="25×35=875"; new="57×60=3420"},
    @{old="15×22=330"; new="60×81=4860"},
    @{old="38×13=494"; new="43×92=3956"},
    @{old="69×47=3243"; new="64×66=4224"},
    @{old="46×41=1886"; new="99×64=6336"},
    @{old="40×95=3800"; new="40×19=760"},
    @{old="55×26=1430"; new="95×87=8265"},
    @{old="62×31=1922"; new="39×21=819"},
    @{old="42×26=1092"; new="84×55=4620"},
    @{old="27×40=1080"; new="83×37=3071"},
    @{old="41×54=2214"; new="84×35=2940"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
